$wb = $excel.ActiveWorkbook
$wsKate = $wb.Worksheets.Item("Kate_hours")
$wsBen  = $wb.Worksheets.Item("Ben_hours")

# --- Kate_hours: new data rows 2-8 --------------------------------------
# Dates / hours first (no shared-string ordering concerns for numbers)
$wsKate.Range("A2").Value = 44998
$wsKate.Range("B2").Value = 1.5
$wsKate.Range("A3").Value = 45005
$wsKate.Range("B3").Value = 3
$wsKate.Range("A4").Value = 45014
$wsKate.Range("B4").Value = 5
$wsKate.Range("A5").Value = 45015
$wsKate.Range("B5").Value = 4
$wsKate.Range("A6").Value = 45020
$wsKate.Range("B6").Value = 2
$wsKate.Range("A7").Value = 45027
$wsKate.Range("B7").Value = 4
$wsKate.Range("A8").Value = 45028
$wsKate.Range("B8").Value = 5

# Task text (C column) + the new "total" label, entered in the exact
# order needed so new shared-string entries come out in the right order.
$wsKate.Range("C6").Value = "Meet with Ben, Jenn, and Eden"
$wsKate.Range("A29").Value = "total"
$wsKate.Range("C4").Value = "Ben and Kate meeting to discuss workflow and assign tasks, OBIS download + initial cleaning"
$wsKate.Range("C5").Value = "GBIF setup and prepping species lists"
$wsKate.Range("C7").Value = "Met with Ben, updated OBIS for ecoregions"
$wsKate.Range("C8").Value = "Fixed loop, top 500 taxize"

# Row 2/3 re-use existing shared strings already present in the workbook.
$wsKate.Range("C2").Value = "Ben and Kate present/discuss eDNA cleaning workflows to Eden"
$wsKate.Range("C3").Value = "discuss Eden's workflow and deliverables of contract"

# Totals formula
$wsKate.Range("B29").Formula = "=SUM(B2:B27)"

# Date formatting: rows 2-6 use m/d/yy (numFmtId 14), rows 7-8 use d-mmm (numFmtId 16)
$wsKate.Range("A2:A6").NumberFormat = "m/d/yy"
$wsKate.Range("A7:A8").NumberFormat = "d-mmm"

# Column widths
$wsKate.Columns.Item(1).ColumnWidth = 11
$wsKate.Columns.Item(2).ColumnWidth = 10.333333333333334

$wsBen.Columns.Item(1).ColumnWidth = 19
$wsBen.Columns.Item(3).ColumnWidth = 73.33333333333333

# --- Selections ----------------------------------------------------------
# Ben_hours: selection becomes A2:XFD5 (not the active sheet afterwards)
[void]$wsBen.Range("A2:XFD5").Select()

# Kate_hours becomes the active sheet/tab with C11 selected
[void]$wsKate.Range("C11").Select()
